$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: merge the three runs (with the spell-check proofErr
#    wrapper around "Wenting") into a single run carrying the full text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Computer Vision week 7 " + [char]8211 + " Wenting Duan " + [char]8211 + " 21/04/2022 " + [char]8211 + " Object Tracking",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Computer Vision week 7 " + [char]8211 + " Wenting Duan " + [char]8211 + " 21/04/2022 " + [char]8211 + " Object Tracking",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Gausian noise" paragraph: fix the typo/spacing and drop the inline
#    spell/grammar-check markers. (Gausian -> Gaussian, "-  a" -> "- a")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Gausian noise-  a statistical",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Gaussian noise- a statistical",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. New content appended after the "Gaussian noise" paragraph:
#    a bold "CVMM:: Pros and Cons" heading followed by three plain
#    paragraphs and a final plain paragraph about Constant Acceleration.
# ---------------------------------------------------------------------------
$pGauss = $d.Paragraphs.Item($d.Paragraphs.Count)
$ins = $pGauss.Range
$ins.Collapse(0)
$ins.InsertParagraphAfter()
$ins.Collapse(0)
$pHeading = $d.Paragraphs.Item($d.Paragraphs.Count)
$pHeading.Range.InsertAfter("CVMM:: Pros and Cons")

$ins = $pHeading.Range
$ins.Collapse(0)
$ins.InsertParagraphAfter()
$ins.Collapse(0)
$pPros = $d.Paragraphs.Item($d.Paragraphs.Count)
$pPros.Range.InsertAfter("This is a popular model for nearly-constant (piecewise) rectilinear motions")

$ins = $pPros.Range
$ins.Collapse(0)
$ins.InsertParagraphAfter()
$ins.Collapse(0)
$pOcclusion = $d.Paragraphs.Item($d.Paragraphs.Count)
$pOcclusion.Range.InsertAfter("Thanks to the velocity components, it can deal with small occlusions.")

$ins = $pOcclusion.Range
$ins.Collapse(0)
$ins.InsertParagraphAfter()
$ins.Collapse(0)
$pTuning = $d.Paragraphs.Item($d.Paragraphs.Count)
$pTuning.Range.InsertAfter("With a careful tuning of the parameters, it can be used for non-rectilinear trajectories, if the latter can be locally approximated to rectilinear ones.")

$ins = $pTuning.Range
$ins.Collapse(0)
$ins.InsertParagraphAfter()
$ins.Collapse(0)
$pAccel = $d.Paragraphs.Item($d.Paragraphs.Count)
$pAccel.Range.InsertAfter("Constant Acceleration & Constant Turn-rate Models")

# Bold the heading text only (not the paragraph mark of the following
# paragraphs) so the bold formatting does not bleed into later paragraphs.
$pHeading = $d.Paragraphs.Item(18)
$rHeading = $pHeading.Range
$rHeading.Font.Bold = 1
$rHeading.Font.BoldBi = 1
